# Edit script for N=10 Graphs.xlsx
# Commit: "Updates to data, MFQ fix, and Graph creation"
#
# 1) Update the N_TT=8 group's AvgWait/TT/Runtime values (columns AK, AL, AM)
#    for rows 2-102 to reflect the corrected/re-run simulation data (MFQ fix).
# 2) Add a new "AVG" summary row (row 103) with AVERAGE() formulas for every
#    data column.
# 3) Nudge the two existing chart objects down by one row so they continue to
#    sit directly below the (now one-row-taller) data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Updated AK/AL/AM values (N_TT group 8: avg wait / TT / Runtime) for rows
#    2 through 102.
# ---------------------------------------------------------------------------
$data = New-Object 'object[,]' 101,3
$data[0,0] = 31.2
$data[0,1] = 5.1963900000000001
$data[0,2] = 417
$data[1,0] = 33.299999999999997
$data[1,1] = 6.05044
$data[1,2] = 323
$data[2,0] = 39.700000000000003
$data[2,1] = 6.3765900000000002
$data[2,2] = 378
$data[3,0] = 33.6
$data[3,1] = 6.92333
$data[3,2] = 403
$data[4,0] = 26.7
$data[4,1] = 8.0909499999999994
$data[4,2] = 249
$data[5,0] = 28.5
$data[5,1] = 4.6902400000000002
$data[5,2] = 467
$data[6,0] = 35.700000000000003
$data[6,1] = 5.6147200000000002
$data[6,2] = 454
$data[7,0] = 22
$data[7,1] = 6.6326200000000002
$data[7,2] = 320
$data[8,0] = 21.9
$data[8,1] = 7.5321400000000001
$data[8,2] = 394
$data[9,0] = 22.4
$data[9,1] = 5.6950000000000003
$data[9,2] = 281
$data[10,0] = 45
$data[10,1] = 5.9410699999999999
$data[10,2] = 722
$data[11,0] = 33.1
$data[11,1] = 5.2497600000000002
$data[11,2] = 436
$data[12,0] = 33.9
$data[12,1] = 6.1201600000000003
$data[12,2] = 394
$data[13,0] = 29.1
$data[13,1] = 6.6916700000000002
$data[13,2] = 199
$data[14,0] = 37.1
$data[14,1] = 10.5783
$data[14,2] = 385
$data[15,0] = 38.200000000000003
$data[15,1] = 9.9502799999999993
$data[15,2] = 484
$data[16,0] = 30.6
$data[16,1] = 4.9667500000000002
$data[16,2] = 426
$data[17,0] = 37.700000000000003
$data[17,1] = 8.1537299999999995
$data[17,2] = 446
$data[18,0] = 48.9
$data[18,1] = 7.3478199999999996
$data[18,2] = 545
$data[19,0] = 48.3
$data[19,1] = 8.0571800000000007
$data[19,2] = 436
$data[20,0] = 33.1
$data[20,1] = 9.0684500000000003
$data[20,2] = 362
$data[21,0] = 23.3
$data[21,1] = 6.5860700000000003
$data[21,2] = 293
$data[22,0] = 39
$data[22,1] = 7.0331000000000001
$data[22,2] = 520
$data[23,0] = 29.5
$data[23,1] = 10.450799999999999
$data[23,2] = 335
$data[24,0] = 32.799999999999997
$data[24,1] = 8.3434500000000007
$data[24,2] = 447
$data[25,0] = 36.299999999999997
$data[25,1] = 5.26607
$data[25,2] = 659
$data[26,0] = 24.9
$data[26,1] = 7.95444
$data[26,2] = 536
$data[27,0] = 53.8
$data[27,1] = 6.4983300000000002
$data[27,2] = 516
$data[28,0] = 28.1
$data[28,1] = 5.0869799999999996
$data[28,2] = 468
$data[29,0] = 22.7
$data[29,1] = 7.8616700000000002
$data[29,2] = 325
$data[30,0] = 36.799999999999997
$data[30,1] = 6.8992899999999997
$data[30,2] = 479
$data[31,0] = 32.1
$data[31,1] = 10.061199999999999
$data[31,2] = 230
$data[32,0] = 49.8
$data[32,1] = 6.50115
$data[32,2] = 593
$data[33,0] = 27.9
$data[33,1] = 6.7602799999999998
$data[33,2] = 319
$data[34,0] = 22.8
$data[34,1] = 5.0650000000000004
$data[34,2] = 633
$data[35,0] = 31.9
$data[35,1] = 6.085
$data[35,2] = 4831
$data[36,0] = 24.2
$data[36,1] = 5.7175399999999996
$data[36,2] = 361
$data[37,0] = 32.700000000000003
$data[37,1] = 9.6199999999999992
$data[37,2] = 350
$data[38,0] = 43.8
$data[38,1] = 7.6003999999999996
$data[38,2] = 372
$data[39,0] = 23.8
$data[39,1] = 5.4971800000000002
$data[39,2] = 581
$data[40,0] = 33.200000000000003
$data[40,1] = 9.0327800000000007
$data[40,2] = 227
$data[41,0] = 20.100000000000001
$data[41,1] = 4.2816700000000001
$data[41,2] = 320
$data[42,0] = 33.5
$data[42,1] = 5.6983300000000003
$data[42,2] = 448
$data[43,0] = 27
$data[43,1] = 7.4517499999999997
$data[43,2] = 210
$data[44,0] = 38.799999999999997
$data[44,1] = 6.2011900000000004
$data[44,2] = 435
$data[45,0] = 20.3
$data[45,1] = 4.5525399999999996
$data[45,2] = 384
$data[46,0] = 29.4
$data[46,1] = 8.8083299999999998
$data[46,2] = 190
$data[47,0] = 30.6
$data[47,1] = 4.9010699999999998
$data[47,2] = 537
$data[48,0] = 31.7
$data[48,1] = 9.4648400000000006
$data[48,2] = 450
$data[49,0] = 29.3
$data[49,1] = 5.8238099999999999
$data[49,2] = 174
$data[50,0] = 34.1
$data[50,1] = 5.9154799999999996
$data[50,2] = 405
$data[51,0] = 32.200000000000003
$data[51,1] = 6.7175000000000002
$data[51,2] = 352
$data[52,0] = 20.100000000000001
$data[52,1] = 4.8410700000000002
$data[52,2] = 200
$data[53,0] = 31.8
$data[53,1] = 5.3291700000000004
$data[53,2] = 1537
$data[54,0] = 32.6
$data[54,1] = 8.0314300000000003
$data[54,2] = 365
$data[55,0] = 29.1
$data[55,1] = 9.3491700000000009
$data[55,2] = 255
$data[56,0] = 35.9
$data[56,1] = 9.2057500000000001
$data[56,2] = 334
$data[57,0] = 23.2
$data[57,1] = 5.6694399999999998
$data[57,2] = 218
$data[58,0] = 37.5
$data[58,1] = 8.5897199999999998
$data[58,2] = 369
$data[59,0] = 31
$data[59,1] = 5.9082499999999998
$data[59,2] = 262
$data[60,0] = 42.1
$data[60,1] = 5.9469399999999997
$data[60,2] = 477
$data[61,0] = 18.3
$data[61,1] = 6.5666700000000002
$data[61,2] = 217
$data[62,0] = 34.799999999999997
$data[62,1] = 4.9028999999999998
$data[62,2] = 401
$data[63,0] = 26.5
$data[63,1] = 5.1776200000000001
$data[63,2] = 248
$data[64,0] = 15
$data[64,1] = 4.2708300000000001
$data[64,2] = 238
$data[65,0] = 27.7
$data[65,1] = 7.7083300000000001
$data[65,2] = 320
$data[66,0] = 21
$data[66,1] = 6.8635700000000002
$data[66,2] = 200
$data[67,0] = 17.8
$data[67,1] = 6.75143
$data[67,2] = 227
$data[68,0] = 23.5
$data[68,1] = 5.6183300000000003
$data[68,2] = 419
$data[69,0] = 20.7
$data[69,1] = 5.7698799999999997
$data[69,2] = 412
$data[70,0] = 32.200000000000003
$data[70,1] = 4.9414300000000004
$data[70,2] = 347
$data[71,0] = 28.4
$data[71,1] = 6.9955600000000002
$data[71,2] = 380
$data[72,0] = 32.5
$data[72,1] = 9.4279799999999998
$data[72,2] = 368
$data[73,0] = 34.4
$data[73,1] = 6.2187299999999999
$data[73,2] = 520
$data[74,0] = 19.5
$data[74,1] = 5.9084099999999999
$data[74,2] = 225
$data[75,0] = 40.200000000000003
$data[75,1] = 5.7540500000000003
$data[75,2] = 484
$data[76,0] = 40.5
$data[76,1] = 5.4641700000000002
$data[76,2] = 349
$data[77,0] = 28.1
$data[77,1] = 5.3377800000000004
$data[77,2] = 388
$data[78,0] = 35
$data[78,1] = 8.0648800000000005
$data[78,2] = 324
$data[79,0] = 43.7
$data[79,1] = 5.6943700000000002
$data[79,2] = 743
$data[80,0] = 35.299999999999997
$data[80,1] = 6.6560300000000003
$data[80,2] = 403
$data[81,0] = 19.7
$data[81,1] = 4.3283300000000002
$data[81,2] = 432
$data[82,0] = 19.600000000000001
$data[82,1] = 4.9192099999999996
$data[82,2] = 202
$data[83,0] = 36.1
$data[83,1] = 5.9096000000000002
$data[83,2] = 649
$data[84,0] = 28.1
$data[84,1] = 6.4541700000000004
$data[84,2] = 259
$data[85,0] = 20.5
$data[85,1] = 4.3748800000000001
$data[85,2] = 443
$data[86,0] = 30.4
$data[86,1] = 5.4055999999999997
$data[86,2] = 737
$data[87,0] = 18.3
$data[87,1] = 5.7472200000000004
$data[87,2] = 448
$data[88,0] = 20.9
$data[88,1] = 4.4325000000000001
$data[88,2] = 232
$data[89,0] = 35.9
$data[89,1] = 8.1281700000000008
$data[89,2] = 423
$data[90,0] = 27.8
$data[90,1] = 6.36972
$data[90,2] = 252
$data[91,0] = 35
$data[91,1] = 8.0448400000000007
$data[91,2] = 500
$data[92,0] = 39.1
$data[92,1] = 5.9371
$data[92,2] = 434
$data[93,0] = 44.1
$data[93,1] = 7.4776199999999999
$data[93,2] = 566
$data[94,0] = 31.2
$data[94,1] = 5.7110300000000001
$data[94,2] = 349
$data[95,0] = 26.3
$data[95,1] = 6.28024
$data[95,2] = 457
$data[96,0] = 30.2
$data[96,1] = 9.0514299999999999
$data[96,2] = 229
$data[97,0] = 22.2
$data[97,1] = 4.7957099999999997
$data[97,2] = 299
$data[98,0] = 20.5
$data[98,1] = 4.3748800000000001
$data[98,2] = 443
$data[99,0] = 30.4
$data[99,1] = 5.4055999999999997
$data[99,2] = 737
$data[100,0] = 18.3
$data[100,1] = 5.7472200000000004
$data[100,2] = 448
$ws.Range("AK2:AM102").Value2 = $data

# ---------------------------------------------------------------------------
# 2) New AVG row (row 103): label in column A, AVERAGE formula across each
#    populated data column (groups separated by blank spacer columns).
# ---------------------------------------------------------------------------
$ws.Range("A103").Value = "AVG"

$ws.Range("B103").Formula = "=AVERAGE(B2:B102)"

$avgCols = @("C","D","F","G","H","I","K","L","M","N","P","Q","R","S","U","V","W","X","Z","AA","AB","AC","AE","AF","AG","AH","AJ","AK","AL","AM")
foreach ($col in $avgCols) {
    $ws.Range($col + "103").Formula = "=AVERAGE(" + $col + "2:" + $col + "102)"
}

# ---------------------------------------------------------------------------
# 3) Move both chart objects down by one default row height (15pt) so they
#    stay anchored below the data table now that row 103 was added.
# ---------------------------------------------------------------------------
$co1 = $ws.ChartObjects().Item(1)
$co2 = $ws.ChartObjects().Item(2)
$co1.Top = $co1.Top + 15
$co2.Top = $co2.Top + 15
